# Generate Report for Handoff
# Adds two new handed-off files (57e7464f-... and 76257ed7-...) to the
# localization status report: one new row in "Overview", and one new row
# per language sheet ("zh-cn", "de-de").

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2          # xlUnderlineStyleSingle
$hyperlinkColor = 15570276       # BGR encoding of RGB FF6495ED (the workbook's custom HyperLink font color)

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.md"
$wsOverview.Range("B6").Value = "e2e\57e7464f-7052-4509-b5fa-2d4ffc83dec5.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = ""
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-13 10:51:04"

$wsOverview.Range("A7").Value = "76257ed7-96db-401e-b597-87a778df0b2f.md"
$wsOverview.Range("B7").Value = "e2e\76257ed7-96db-401e-b597-87a778df0b2f.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = ""
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-13 10:51:04"

$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/57e7464f70524509b5fa2d4ffc83dec557e7464f/e2e/57e7464f-7052-4509-b5fa-2d4ffc83dec5.md", "", "", "e2e\57e7464f-7052-4509-b5fa-2d4ffc83dec5.md")
$wsOverview.Range("B6").Font.Underline = $hyperlinkUnderline
$wsOverview.Range("B6").Font.Color = $hyperlinkColor

$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/76257ed796db401eb59787a778df0b2f76257ed/e2e/76257ed7-96db-401e-b597-87a778df0b2f.md", "", "", "e2e\76257ed7-96db-401e-b597-87a778df0b2f.md")
$wsOverview.Range("B7").Font.Underline = $hyperlinkUnderline
$wsOverview.Range("B7").Font.Color = $hyperlinkColor

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A6").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.93a72b73dff4fbf0545eafe0775adbb35b50061a.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-08-13 10:50:54"
$wsZhCn.Range("H6").NumberFormat = $dateFormat
$wsZhCn.Range("I6").Value = ""
$wsZhCn.Range("J6").Value = ""
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K6").NumberFormat = $dateFormat
$wsZhCn.Range("L6").Value = ""
$wsZhCn.Range("M6").Value = "True"
$wsZhCn.Range("N6").Value = ""
$wsZhCn.Range("O6").Value = "False"
$wsZhCn.Range("P6").Value = ""

$wsZhCn.Range("A7").Value = "76257ed7-96db-401e-b597-87a778df0b2f.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "76257ed7-96db-401e-b597-87a778df0b2f.9759c485f51b97edfb63a787712e65eefe3859e2.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-13 10:50:54"
$wsZhCn.Range("H7").NumberFormat = $dateFormat
$wsZhCn.Range("I7").Value = ""
$wsZhCn.Range("J7").Value = ""
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K7").NumberFormat = $dateFormat
$wsZhCn.Range("L7").Value = ""
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("N7").Value = ""
$wsZhCn.Range("O7").Value = "False"
$wsZhCn.Range("P7").Value = ""

$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/57e7464f70524509b5fa2d4ffc83dec557e7464f/e2e/57e7464f-7052-4509-b5fa-2d4ffc83dec5.md", "", "", "57e7464f-7052-4509-b5fa-2d4ffc83dec5.md")
$wsZhCn.Range("A6").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("A6").Font.Color = $hyperlinkColor

$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/76257ed796db401eb59787a778df0b2f76257ed/e2e/76257ed7-96db-401e-b597-87a778df0b2f.md", "", "", "76257ed7-96db-401e-b597-87a778df0b2f.md")
$wsZhCn.Range("A7").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("A7").Font.Color = $hyperlinkColor

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A6").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.93a72b73dff4fbf0545eafe0775adbb35b50061a.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-08-13 10:51:04"
$wsDeDe.Range("H6").NumberFormat = $dateFormat
$wsDeDe.Range("I6").Value = ""
$wsDeDe.Range("J6").Value = ""
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K6").NumberFormat = $dateFormat
$wsDeDe.Range("L6").Value = ""
$wsDeDe.Range("M6").Value = "True"
$wsDeDe.Range("N6").Value = ""
$wsDeDe.Range("O6").Value = "False"
$wsDeDe.Range("P6").Value = ""

$wsDeDe.Range("A7").Value = "76257ed7-96db-401e-b597-87a778df0b2f.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "76257ed7-96db-401e-b597-87a778df0b2f.9759c485f51b97edfb63a787712e65eefe3859e2.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-13 10:51:04"
$wsDeDe.Range("H7").NumberFormat = $dateFormat
$wsDeDe.Range("I7").Value = ""
$wsDeDe.Range("J7").Value = ""
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K7").NumberFormat = $dateFormat
$wsDeDe.Range("L7").Value = ""
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("N7").Value = ""
$wsDeDe.Range("O7").Value = "False"
$wsDeDe.Range("P7").Value = ""

$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/57e7464f70524509b5fa2d4ffc83dec557e7464f/e2e/57e7464f-7052-4509-b5fa-2d4ffc83dec5.md", "", "", "57e7464f-7052-4509-b5fa-2d4ffc83dec5.md")
$wsDeDe.Range("A6").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("A6").Font.Color = $hyperlinkColor

$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/76257ed796db401eb59787a778df0b2f76257ed/e2e/76257ed7-96db-401e-b597-87a778df0b2f.md", "", "", "76257ed7-96db-401e-b597-87a778df0b2f.md")
$wsDeDe.Range("A7").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("A7").Font.Color = $hyperlinkColor

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P7"))

Write-Output "Report updated for handback"
